# Updating score criteria for FL frameworks.
#
# The original workbook has a "T.6" sheet comparing FL frameworks across
# 5 criteria (Ease of use / Flexibility and customisability / Real-world
# applicability), ranked 1st..5th.
#
# This edit:
#   1. Duplicates "T.6" -> "T.6 (2)" (a working copy) BEFORE changing "T.6",
#      since the duplicate starts from the pre-edit layout.
#   2. Reworks "T.6" into the final scoring table: adds a "Developer" column,
#      re-orders FATE/TFF, and collapses the score columns down to three
#      named criteria (Setup and Configuration / Examples and Tutorials /
#      Custom Algorithm Implementation / Adaptability to Various Use Cases /
#      Industry Adoption) plus a hard-coded "Average Score", with a note
#      about methodology in column O.
#   3. Reworks "T.6 (2)" into a smaller, formula-driven variant (3 score
#      columns averaged live via SUM(..)/3) headed "Close to Real-world
#      Settings??", which ends up the active tab.

$wb = $excel.ActiveWorkbook
$ws6 = $wb.Worksheets.Item("T.6")

# ---------------------------------------------------------------------
# Step 1: duplicate "T.6" (in its original, unedited state) to the end of
# the workbook. Excel names the copy "T.6 (2)" automatically.
# ---------------------------------------------------------------------
$ws6.Copy($null, $ws6)
$ws6b = $wb.Worksheets.Item("T.6 (2)")

# =======================================================================
# Step 2: rework "T.6" into its final layout
# =======================================================================

# Remove the old merges so every cell can be written independently.
$ws6.Cells.UnMerge()

# --- Row 4: grouped headers -------------------------------------------
$ws6.Range("D4:N20").ClearContents()

$ws6.Range("G4").Value = "Ease of use"
$ws6.Range("I4").Value = "Flexibility and customisability"
$ws6.Range("K4").Value = "Real-world applicability"
$ws6.Range("O4").Value = "Look and see what other people have done. To evaluate FL frameworks. Look if there is an standard. Explain why did you picked. Also look at weights."

$ws6.Range("G4:H4").Merge()
$ws6.Range("I4:J4").Merge()
$ws6.Range("D4:F4").Merge()

# --- Row 5: column headers ---------------------------------------------
$ws6.Range("D5").Value = "Rank"
$ws6.Range("E5").Value = "Framework"
$ws6.Range("F5").Value = "Developer"
$ws6.Range("G5").Value = "Setup and Configuration"
$ws6.Range("H5").Value = "Examples and Tutorials"
$ws6.Range("I5").Value = "Custom Algorithm Implementation"
$ws6.Range("J5").Value = "Adaptability to Various Use Cases"
$ws6.Range("K5").Value = "Industry Adoption"
$ws6.Range("L5").Value = "Average Score"

# --- Rows 6-10: data -----------------------------------------------------
$ws6.Range("D6").Value = "1st"
$ws6.Range("E6").Value = "FedML"
$ws6.Range("F6").Value = "FedML AI Inc."
$ws6.Range("G6").Value = 8
$ws6.Range("H6").Value = 8
$ws6.Range("I6").Value = 7
$ws6.Range("J6").Value = 8
$ws6.Range("K6").Value = 7
$ws6.Range("L6").Value = 7.6

$ws6.Range("D7").Value = "2nd"
$ws6.Range("E7").Value = "Flower"
$ws6.Range("F7").Value = "Adap GmbH"
$ws6.Range("G7").Value = 7
$ws6.Range("H7").Value = 8
$ws6.Range("I7").Value = 8
$ws6.Range("J7").Value = 8
$ws6.Range("K7").Value = 6
$ws6.Range("L7").Value = 7.4

$ws6.Range("D8").Value = "3rd"
$ws6.Range("E8").Value = "FATE"
$ws6.Range("F8").Value = "WeBank & Linux Foundation"
$ws6.Range("G8").Value = 6
$ws6.Range("H8").Value = 7
$ws6.Range("I8").Value = 7
$ws6.Range("J8").Value = 7
$ws6.Range("K8").Value = 7
$ws6.Range("L8").Value = 6.8

$ws6.Range("D9").Value = "4th"
$ws6.Range("E9").Value = "TFF"
$ws6.Range("F9").Value = "Google"
$ws6.Range("G9").Value = 6
$ws6.Range("H9").Value = 6
$ws6.Range("I9").Value = 9
$ws6.Range("J9").Value = 6
$ws6.Range("K9").Value = 6
$ws6.Range("L9").Value = 6.6

$ws6.Range("D10").Value = "5th"
$ws6.Range("E10").Value = "PySyft"
$ws6.Range("F10").Value = "OpenMined"
$ws6.Range("G10").Value = 5
$ws6.Range("H10").Value = 5
$ws6.Range("I10").Value = 6
$ws6.Range("J10").Value = 5
$ws6.Range("K10").Value = 4
$ws6.Range("L10").Value = 5

# --- Formatting ---------------------------------------------------------
$ws6.Rows.Item(4).RowHeight = 28.8
$ws6.Rows.Item(5).RowHeight = 43.2

$headerRange = $ws6.Range("D4:L5")
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108
$headerRange.WrapText = $true
$headerRange.Font.Name = "Calibri"
$headerRange.Font.Bold = $true
$ws6.Range("D4:F4").Font.Bold = $false

$dataBody = $ws6.Range("D4:L10")
$dataBody.Borders.LineStyle = 1

$dataValues = $ws6.Range("D6:L10")
$dataValues.HorizontalAlignment = -4108
$dataValues.VerticalAlignment = -4108
$dataValues.WrapText = $true
$dataValues.Font.Name = "Calibri"
$dataValues.Font.Bold = $false
$ws6.Range("D6:D10,F6:K10").Font.Bold = $false
$ws6.Range("E6:E10,L6:L10").Font.Bold = $true

$ws6.Range("O4").Font.Name = "Aptos Narrow"
$ws6.Range("O4").Font.Bold = $true

$ws6.Columns.Item("E").ColumnWidth = 15.33203125
$ws6.Columns.Item("F").ColumnWidth = 15.33203125
$ws6.Columns.Item("G").ColumnWidth = 13.109375
$ws6.Columns.Item("I").ColumnWidth = 15.88671875
$ws6.Columns.Item("J").ColumnWidth = 19.33203125
$ws6.Columns.Item("K").ColumnWidth = 16.6640625
$ws6.Columns.Item("O").ColumnWidth = 130.77734375

# Leftover blank-but-formatted placeholder block below the table
# (rows 13-20, columns D:N) - carries no text, only touches the used range.
$ws6.Range("D13:N20").HorizontalAlignment = -4108
$ws6.Range("D13:N20").VerticalAlignment = -4108
$ws6.Range("D13:N20").WrapText = $true
$ws6.Range("E14:K14,E17:K17").Font.Name = "Calibri"
$ws6.Range("E14,E17").Font.Bold = $true
$ws6.Range("F14:K14,F17:K17").Font.Bold = $false

[void]$ws6.Range("D4:L10").Select()

# =======================================================================
# Step 3: rework "T.6 (2)" (the duplicate) into the formula-driven variant
# =======================================================================
$ws6b.Cells.UnMerge()
$ws6b.Range("D4:N20").ClearContents()

$ws6b.Range("G4").Value = "Ease of use"
$ws6b.Range("I4").Value = "Close to Real-world Settings??"
$ws6b.Range("N4").Value = "Look and see what other people have done. To evaluate FL frameworks. Look if there is an standard. Explain why did you picked. Also look at weights."

$ws6b.Range("G4:H4").Merge()
$ws6b.Range("D4:F4").Merge()

$ws6b.Range("D5").Value = "Rank"
$ws6b.Range("E5").Value = "Framework"
$ws6b.Range("F5").Value = "Developer"
$ws6b.Range("G5").Value = " Setup and Configuration"
$ws6b.Range("H5").Value = "Adaptability to Various Use Cases"
$ws6b.Range("I5").Value = "Examples and Tutorials"
$ws6b.Range("J5").Value = "Average"

$ws6b.Range("D6").Value = "1st"
$ws6b.Range("E6").Value = "FedML"
$ws6b.Range("F6").Value = "FedML AI Inc."
$ws6b.Range("G6").Value = 7
$ws6b.Range("H6").Value = 8
$ws6b.Range("I6").Value = 8
$ws6b.Range("J6").Formula = "=SUM(G6:I6)/3"

$ws6b.Range("D7").Value = "2nd"
$ws6b.Range("E7").Value = "Flower"
$ws6b.Range("F7").Value = "Adap GmbH"
$ws6b.Range("G7").Value = 7
$ws6b.Range("H7").Value = 8
$ws6b.Range("I7").Value = 7
$ws6b.Range("J7").Formula = "=SUM(G7:I7)/3"

$ws6b.Range("D8").Value = "3rd"
$ws6b.Range("E8").Value = "FATE"
$ws6b.Range("F8").Value = "WeBank & Linux Foundation"
$ws6b.Range("G8").Value = 6
$ws6b.Range("H8").Value = 7
$ws6b.Range("I8").Value = 7
$ws6b.Range("J8").Formula = "=SUM(G8:I8)/3"

$ws6b.Range("D9").Value = "4th"
$ws6b.Range("E9").Value = "TFF"
$ws6b.Range("F9").Value = "Google"
$ws6b.Range("G9").Value = 7
$ws6b.Range("H9").Value = 6
$ws6b.Range("I9").Value = 6
$ws6b.Range("J9").Formula = "=SUM(G9:I9)/3"

$ws6b.Range("D10").Value = "5th"
$ws6b.Range("E10").Value = "PySyft"
$ws6b.Range("F10").Value = "OpenMined"
$ws6b.Range("G10").Value = 5
$ws6b.Range("H10").Value = 5
$ws6b.Range("I10").Value = 5
$ws6b.Range("J10").Formula = "=SUM(G10:I10)/3"

# --- Formatting ---------------------------------------------------------
$ws6b.Rows.Item(4).RowHeight = 28.8
$ws6b.Rows.Item(5).RowHeight = 28.8

$headerRange2 = $ws6b.Range("D4:J5")
$headerRange2.HorizontalAlignment = -4108
$headerRange2.VerticalAlignment = -4108
$headerRange2.WrapText = $true
$headerRange2.Font.Name = "Calibri"
$headerRange2.Font.Bold = $true
$ws6b.Range("D4:F4").Font.Bold = $false

$dataBody2 = $ws6b.Range("D4:J10")
$dataBody2.Borders.LineStyle = 1

$dataValues2 = $ws6b.Range("D6:J10")
$dataValues2.HorizontalAlignment = -4108
$dataValues2.VerticalAlignment = -4108
$dataValues2.WrapText = $true
$dataValues2.Font.Name = "Calibri"
$dataValues2.Font.Bold = $false
$ws6b.Range("E6:E10,J6:J10").Font.Bold = $true

$ws6b.Range("N4").Font.Name = "Aptos Narrow"
$ws6b.Range("N4").Font.Bold = $true

$ws6b.Columns.Item("E").ColumnWidth = 15.33203125
$ws6b.Columns.Item("F").ColumnWidth = 15.33203125
$ws6b.Columns.Item("G").ColumnWidth = 13.109375
$ws6b.Columns.Item("H").ColumnWidth = 18.109375
$ws6b.Columns.Item("I").ColumnWidth = 15.44140625
$ws6b.Columns.Item("J").ColumnWidth = 15.44140625
$ws6b.Columns.Item("K").ColumnWidth = 15.44140625
$ws6b.Columns.Item("N").ColumnWidth = 130.77734375

# Leftover blank-but-formatted placeholder block below the table
# (rows 13-20, columns D:M) - carries no text, only touches the used range.
$ws6b.Range("D13:M20").HorizontalAlignment = -4108
$ws6b.Range("D13:M20").VerticalAlignment = -4108
$ws6b.Range("D13:M20").WrapText = $true
$ws6b.Range("E14:K14,E17:K17").Font.Name = "Calibri"
$ws6b.Range("E14,E17").Font.Bold = $true
$ws6b.Range("F14:K14,F17:K17").Font.Bold = $false

# Final state: "T.6 (2)" is the active / selected tab.
[void]$ws6b.Activate()
[void]$ws6b.Range("J29").Select()
